$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# First table (column A, rows 2-6): append " (Iago)" to the tasks assigned to Iago
$ws.Range("A2").Value = "Debugs (Iago)"
$ws.Range("A3").Value = "Correct bugs (Iago)"
$ws.Range("A6").Value = "Make java program to encapsulate tasks (Iago)"

# Second table (column D, rows 10-14): same tasks listed again, mirror the edits
$ws.Range("D10").Value = "Debugs (Iago)"
$ws.Range("D11").Value = "Correct bugs (Iago)"
$ws.Range("D14").Value = "Make java program to encapsulate tasks (Iago)"

# Reflect the active cell being left on D14 after the last edit
$ws.Range("D14").Select()
